# Applies the changes described in the diff:
#  - Row 16: set H16 = 5 (style already s="2"); add I16 = 5 (style s="8"); add J16 = 5 (style s="9")
#  - Row 28: set C28 = 5; clear D28/E28/F28 (keep style s="2"); set H28 = 5 (style already s="2"); add I28 = 5 (style s="8")
#  - Update the active selection on the sheet to C25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 updates ---
# H16 already carries style s="2"; just set its value.
$ws.Range("H16").Value = 5

# I16 needs style s="8" (same style as I22/I25). Copy format then set value.
$ws.Range("I22").Copy()
$ws.Range("I16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I16").Value = 5

# J16 needs style s="9" (same style as J25). Copy format then set value.
$ws.Range("J25").Copy()
$ws.Range("J16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J16").Value = 5

$excel.CutCopyMode = 0

# --- Row 28 updates ---
$ws.Range("C28").Value = 5
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("H28").Value = 5

# I28 needs style s="8" (same style as I22/I25). Copy format then set value.
$ws.Range("I22").Copy()
$ws.Range("I28").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I28").Value = 5

$excel.CutCopyMode = 0

# --- Update selection / active cell ---
$ws.Range("C25").Select()
